$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Correção do campo Data em acervos tridimensionais":
# remove the obsolete "Data" column (E) from the header row; the
# columns to its right shift left to fill the gap.
$ws.Range("E1").EntireColumn.Delete() | Out-Null
$ws.Range("H10").Select() | Out-Null
